$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '60.688.14'
$ws.Range("E2").Value = '  -2.33%  '
$ws.Range("D3").Value = '2.366.58'
$ws.Range("E3").Value = '  -3.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.82'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -1.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.48'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -2.62%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("D9").Value = '2.366.49'
$ws.Range("E9").Value = '  -3.08%  '
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.08'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -2.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.340'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.81'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("E15").Value = '  -1.53%  '
$ws.Range("D16").Value = '2.824.35'
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("D17").Value = '60.708.04'
$ws.Range("E17").Value = '  -2.21%  '
$ws.Range("D18").Value = '2.367.36'
$ws.Range("E18").Value = '  -2.26%  '
$ws.Range("E19").Value = '  -2.72%  '
$ws.Range("E20").Value = '  +1.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.06'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -2.13%  '
$ws.Range("E22").Value = '  -1.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.08'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +1.75%  '
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("E25").Value = '  -6.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.62'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -1.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.58'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -8.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '571.10'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -6.05%  '
$ws.Range("D30").Value = '0.0₃0911'
$ws.Range("E30").Value = '  -3.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.85'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -1.58%  '
$ws.Range("E32").Value = '  -5.46%  '
$ws.Range("E33").Value = '  -2.46%  '
$ws.Range("E34").Value = '  -5.98%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.61'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -5.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.368'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -2.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.38'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -3.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '146.87'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -1.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.11'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("E41").Value = '  -4.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  -3.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.87'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -3.92%  '
$ws.Range("E45").Value = '  -4.47%  '
$ws.Range("E46").Value = '  +19.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '139.46'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -2.11%  '
$ws.Range("E48").Value = '  -3.47%  '
$ws.Range("E49").Value = '  -3.34%  '
$ws.Range("E50").Value = '  -3.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.28'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -0.73%  '
